$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("B2").Value = 0.1868852459016394
$ws.Range("C2").Value = 0.580327868852459
$ws.Range("J2").Value = 0.01967213114754099
$ws.Range("O2").Value = 0.003278688524590164
$ws.Range("P2").Value = 0.1311475409836066
$ws.Range("S2").Value = 0.07868852459016394

# Row 3
$ws.Range("C3").Value = 0.02702702702702703
$ws.Range("J3").Value = 0.02702702702702703
$ws.Range("P3").Value = 0.7243243243243244
$ws.Range("S3").Value = 0.2216216216216216

# Row 4
$ws.Range("J4").Value = 0.04761904761904762
$ws.Range("P4").Value = 0.7380952380952381
$ws.Range("S4").Value = 0.2142857142857143

# Row 6
$ws.Range("B6").Value = 0.05581395348837209
$ws.Range("D6").Value = 0.004651162790697674
$ws.Range("F6").Value = 0.04651162790697674
$ws.Range("J6").Value = 0.2604651162790698
$ws.Range("O6").Value = 0.01395348837209302
$ws.Range("Q6").Value = 0.2325581395348837
$ws.Range("R6").Value = 0.04651162790697674
$ws.Range("S6").Value = 0.3395348837209302

# Row 7
$ws.Range("B7").Value = 0.1515151515151515
$ws.Range("D7").Value = 0.0101010101010101
$ws.Range("F7").Value = 0.06060606060606061
$ws.Range("J7").Value = 0.1060606060606061
$ws.Range("O7").Value = 0.0101010101010101
$ws.Range("Q7").Value = 0.1919191919191919
$ws.Range("R7").Value = 0.07575757575757576
$ws.Range("S7").Value = 0.3939393939393939

# Row 8
$ws.Range("B8").Value = 0.1155913978494624
$ws.Range("D8").Value = 0.01881720430107527
$ws.Range("F8").Value = 0.0456989247311828
$ws.Range("J8").Value = 0.1129032258064516
$ws.Range("O8").Value = 0.02419354838709677
$ws.Range("Q8").Value = 0.2043010752688172
$ws.Range("R8").Value = 0.1102150537634409
$ws.Range("S8").Value = 0.3682795698924731

# Row 9
$ws.Range("B9").Value = 0.1145374449339207
$ws.Range("D9").Value = 0.02643171806167401
$ws.Range("E9").Value = 0.004405286343612335
$ws.Range("F9").Value = 0.07048458149779736
$ws.Range("J9").Value = 0.09251101321585903
$ws.Range("O9").Value = 0.03083700440528634
$ws.Range("Q9").Value = 0.2026431718061674
$ws.Range("R9").Value = 0.06607929515418502
$ws.Range("S9").Value = 0.3920704845814978

# Row 10
$ws.Range("B10").Value = 0.1072555205047319
$ws.Range("D10").Value = 0.02050473186119874
$ws.Range("E10").Value = 0.0007886435331230284
$ws.Range("F10").Value = 0.0749211356466877
$ws.Range("J10").Value = 0.1167192429022082
$ws.Range("O10").Value = 0.01261829652996845
$ws.Range("Q10").Value = 0.2026813880126183
$ws.Range("R10").Value = 0.09148264984227129
$ws.Range("S10").Value = 0.3730283911671924

# Row 11
$ws.Range("G11").Value = 0.1516129032258065
$ws.Range("J11").Value = 0.08064516129032258
$ws.Range("K11").Value = 0.2032258064516129
$ws.Range("L11").Value = 0.5387096774193548
$ws.Range("S11").Value = 0.02580645161290323

# Row 12
$ws.Range("G12").Value = 0.7078651685393258
$ws.Range("J12").Value = 0.1797752808988764
$ws.Range("K12").Value = 0.01685393258426966
$ws.Range("L12").Value = 0.06741573033707865
$ws.Range("S12").Value = 0.02808988764044944

# Row 13
$ws.Range("G13").Value = 0.5957446808510638
$ws.Range("J13").Value = 0.2978723404255319
$ws.Range("S13").Value = 0.1063829787234043

# Row 14
$ws.Range("G14").Value = 0.6666666666666666
$ws.Range("J14").Value = 0.3333333333333333

# Row 15
$ws.Range("F15").Value = 0.009049773755656109
$ws.Range("H15").Value = 0.1131221719457014
$ws.Range("I15").Value = 0.07692307692307693
$ws.Range("J15").Value = 0.4343891402714932
$ws.Range("K15").Value = 0.07239819004524888
$ws.Range("M15").Value = 0.004524886877828055
$ws.Range("O15").Value = 0.04524886877828054
$ws.Range("S15").Value = 0.2443438914027149

# Row 16
$ws.Range("F16").Value = 0.02030456852791878
$ws.Range("H16").Value = 0.1979695431472081
$ws.Range("I16").Value = 0.1015228426395939
$ws.Range("J16").Value = 0.3959390862944163
$ws.Range("K16").Value = 0.08629441624365482
$ws.Range("M16").Value = 0.03045685279187817
$ws.Range("O16").Value = 0.07106598984771574
$ws.Range("S16").Value = 0.09644670050761421

# Row 17
$ws.Range("F17").Value = 0.02150537634408602
$ws.Range("H17").Value = 0.1505376344086022
$ws.Range("I17").Value = 0.1053763440860215
$ws.Range("J17").Value = 0.3956989247311828
$ws.Range("K17").Value = 0.1053763440860215
$ws.Range("M17").Value = 0.02580645161290323
$ws.Range("N17").Value = 0.006451612903225806
$ws.Range("O17").Value = 0.05591397849462366
$ws.Range("S17").Value = 0.1333333333333333

# Row 18
$ws.Range("F18").Value = 0.02040816326530612
$ws.Range("H18").Value = 0.163265306122449
$ws.Range("I18").Value = 0.1020408163265306
$ws.Range("J18").Value = 0.3979591836734694
$ws.Range("K18").Value = 0.09693877551020408
$ws.Range("M18").Value = 0.02551020408163265
$ws.Range("O18").Value = 0.08673469387755102
$ws.Range("S18").Value = 0.1071428571428571

# Row 19
$ws.Range("F19").Value = 0.01394585726004922
$ws.Range("H19").Value = 0.1706316652994257
$ws.Range("I19").Value = 0.0992616899097621
$ws.Range("J19").Value = 0.3863822805578343
$ws.Range("K19").Value = 0.1156685808039377
$ws.Range("M19").Value = 0.01968826907301066
$ws.Range("N19").Value = 0.002461033634126333
$ws.Range("O19").Value = 0.07465135356849877
$ws.Range("S19").Value = 0.1173092698933552
